$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (GOOG)
$ws.Range("D2").Value = 0.003130648843205768
$ws.Range("E2").Value = 0.01850765251460834
$ws.Range("G2").Value = 0.51
$ws.Range("H2").Value = 22091.34

# Row 3 (MSFT)
$ws.Range("D3").Value = 0.001316074188654831
$ws.Range("E3").Value = 0.01585281793820317

# Row 4 (FB)
$ws.Range("D4").Value = 0.001184645926661158
$ws.Range("E4").Value = 0.02227509086647682

# Row 5 (AAPL)
$ws.Range("D5").Value = 0.001043213982647104
$ws.Range("E5").Value = 0.02021608906402109

# Row 6 (ARKF)
$ws.Range("D6").Value = 0.002176021606322795
$ws.Range("E6").Value = 0.02342341165530148

# Row 7 (ENPH)
$ws.Range("D7").Value = 0.004115310161708273
$ws.Range("E7").Value = 0.04909455604170841
$ws.Range("G7").Value = 0.01
$ws.Range("H7").Value = 433.16

# Row 8 (AMZN)
$ws.Range("D8").Value = 0.00006118006989292327
$ws.Range("E8").Value = 0.01791952477032176

# Row 9 (TSLA)
$ws.Range("D9").Value = 0.005117839680748134
$ws.Range("E9").Value = 0.04199428619745113
$ws.Range("G9").Value = 0.07000000000000001
$ws.Range("H9").Value = 3032.14

# Row 10 (PLTR)
$ws.Range("D10").Value = 0.01008735305826615
$ws.Range("E10").Value = 0.06709706348823545

# Row 11 (RDS-A)
$ws.Range("D11").Value = 0.00362844919629976
$ws.Range("E11").Value = 0.02775403592979557
$ws.Range("G11").Value = 0.27
$ws.Range("H11").Value = 11695.41
